# This commit just re-saved the workbook after clicking a different cell on
# Sheet2 (selection moved from E4 to B5) - no formulas or data were edited.
# Saving also makes Excel recompute the volatile RAND()-based formulas on
# Sheet2 (B2/B3) together with their dependents Sheet1!A2 and Sheet3!A2.
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 is (and stays) the selected/active tab; just move the selection.
$ws2.Activate()
$ws2.Range("B5").Select()

# Recalculate so the volatile formulas get fresh cached values, as happens
# whenever the workbook is recalculated/saved.
$excel.CalculateFull()
